$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Cells.Item(7, 1).Value = 'The Trapping Mechanism at the AlGaN/GaN Interface and the Turn-On Characteristics of the p-GaN Direct-Coupled FET Logic Inverters'
$ws.Cells.Item(7, 2).Value = 2024
$ws.Cells.Item(7, 3).Value = 'Nanomaterials'
$ws.Cells.Item(7, 4).Value = 'Junfeng Yu; Jihong Ding; Tao Wang; Yukai Huang; Wenzhang Du; Jiao Liang; Hongping Ma; Qingchun Zhang; Liang Li; Wei Huang; Wei Zhang'
$ws.Cells.Item(7, 5).Value = '10.3390/nano14241984'
$ws.Cells.Item(7, 6).Value = 'https://doi.org/10.3390/nano14241984'
$ws.Cells.Item(7, 7).Value = 'B3_non_cmos_logic_context; mentions_p_channel; mentions_n_channel; mentions_delay_metrics'
$ws.Cells.Item(7, 8).Value = 'NODE_B2_GaN_CMOS_Logic_Demos'
$ws.Cells.Item(7, 9).Value = '"GaN" complementary logic'
$ws.Cells.Item(7, 10).Value = '2026-02-03T04:29:09.402327+00:00'
$ws.Cells.Item(7, 11).Value = 'The trapping mechanism at the AlGaN/GaN interface in the p-GaN high electron mobility transistors (HEMTs) and its impact on the turn-on characteristics of direct-coupled FET logic (DCFL) inverters were investigated across various supply voltages (VDD) and test frequencies (fm). The frequency-conductance method identified two trap states at the AlGaN/GaN interface (trap activation energy Ec-ET ranges from 0.345 eV to 0.363 eV and 0.438 eV to 0.47 eV). As VDD increased from 1.5 V to 5 V, the interface traps captured more electrons, increasing the channel resistance (Rchannel) and drift-region resistance (Rdrift) of the p-GaN HEMTs and raising the low-level voltage (VOL) from 0.56 V to 1.01 V. At fm = 1 kHz, sufficient trapping and de-trapping led to a delay of 220 µs and a VOL instability of 320 mV. Additionally, as fm increased from 1 kHz to 200 kHz, a positive shift in the threshold voltage of p-GaN HEMTs occurred due to the dominance of trapping. This shift caused VOL to rise from 1.02 V to 1.40 V and extended the fall time (tfall) from 153 ns to 1 µs. This investigation enhances the understanding of DCFL GaN inverters’ behaviors from the perspective of device physics on power switching applications.'
$ws.Cells.Item(7, 12).Value = 'Ambiguous Node-B evidence (B3-only or weak match)'

# Row 8
$ws.Cells.Item(8, 1).Value = 'Vertical integration: a key concept for future flexible and printed electronics'
$ws.Cells.Item(8, 2).Value = 2022
$ws.Cells.Item(8, 4).Value = 'Heesung Han; Chang‐Hyun Kim; Sungyeop Jung'
$ws.Cells.Item(8, 5).Value = '10.1088/2058-8585/ac7929'
$ws.Cells.Item(8, 6).Value = 'https://doi.org/10.1088/2058-8585/ac7929'
$ws.Cells.Item(8, 7).Value = 'B3_non_cmos_logic_context'
$ws.Cells.Item(8, 8).Value = 'OA_CITER'
$ws.Cells.Item(8, 9).Value = '10.1002/aelm.201600046'
$ws.Cells.Item(8, 10).Value = '2026-02-03T04:34:40.343693+00:00'
$ws.Cells.Item(8, 11).Value = 'Abstract This review aims at introducing a vertical integration approach as a promising new driver of field-effect transistor circuits and systems, which can overcome limitations of flexible and printed electronics. A large portion of this article is dedicated to systematic categorization and illustration of important and recent demonstrations of vertically integrated circuit building blocks incorporating organic materials, metal-oxide semiconductors, nanotubes, and nanosheets. Key features of each of these materials platforms and their fabrication issues are also broadly discussed. Then, common technical requirements for high-performance devices and circuits are critically evaluated in view of the specific geometrical and theoretical aspects of vertically connected complementary digital logic inverters. Prominent technological opportunities and challenges for the vertical and three-dimensional transistor integration are finally addressed to further motivate active multidisciplinary research on related materials, devices, and systems.'
$ws.Cells.Item(8, 12).Value = 'Ambiguous Node-B evidence (B3-only or weak match)'
